$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-30 Tuesday" "2024-07-31 Wednesday"

Replace-Text "180×5=" "779×6="
Replace-Text "138×2=" "577×2="
Replace-Text "857×9=" "310×5="
Replace-Text "162×8=" "826×8="
Replace-Text "520×8=" "266×5="

Replace-Text "107×2=" "914×7="
Replace-Text "979×2=" "301×5="
Replace-Text "188×2=" "530×6="
Replace-Text "488×4=" "133×6="
Replace-Text "546×5=" "574×5="

Replace-Text "726×4=" "449×8="
Replace-Text "698×2=" "878×9="
Replace-Text "311×8=" "157×9="
Replace-Text "451×4=" "664×4="
Replace-Text "468×7=" "453×8="

Replace-Text "541×5=" "242×6="
Replace-Text "807×6=" "828×7="
Replace-Text "252×6=" "583×8="
Replace-Text "397×7=" "364×5="
Replace-Text "828×6=" "843×9="

Replace-Text "774×9=" "650×8="
Replace-Text "821×6=" "567×8="
Replace-Text "967×6=" "655×8="
Replace-Text "993×6=" "507×7="
Replace-Text "792×5=" "650×6="
